# Insert a new weekly price record as row 15 (Fecha 2022-08-04), pushing the
# existing rows 15-22 down to 16-23.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(15, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(15, 4).Value = 44777
$ws.Cells.Item(15, 5).Value = 15
$ws.Cells.Item(15, 6).Value = 100112003
$ws.Cells.Item(15, 7).Value = "Ajo"
$ws.Cells.Item(15, 8).Value = "Chino"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 200
$ws.Cells.Item(15, 11).Value = 24000
$ws.Cells.Item(15, 12).Value = 25000
$ws.Cells.Item(15, 13).Value = 24500
$ws.Cells.Item(15, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(15, 15).Value = "China"
$ws.Cells.Item(15, 16).Value = 2450
$ws.Cells.Item(15, 17).Value = 10
$ws.Cells.Item(15, 18).Value = "Hortaliza"
